$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2026-02-18 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2026-02-19 Thursday", 2) | Out-Null

# Update the multiplication problems in the table.
# The table has 20 rows; only rows 1, 5, 10, 15, 20 contain data (5 columns each),
# the other rows are empty spacer rows.
$t = $d.Tables.Item(1)

$replacements = @(
    @{ Row = 1;  Col = 1; New = "606×4=2424" },
    @{ Row = 1;  Col = 2; New = "582×7=4074" },
    @{ Row = 1;  Col = 3; New = "912×4=3648" },
    @{ Row = 1;  Col = 4; New = "921×3=2763" },
    @{ Row = 1;  Col = 5; New = "742×6=4452" },

    @{ Row = 5;  Col = 1; New = "323×8=2584" },
    @{ Row = 5;  Col = 2; New = "559×3=1677" },
    @{ Row = 5;  Col = 3; New = "472×8=3776" },
    @{ Row = 5;  Col = 4; New = "430×6=2580" },
    @{ Row = 5;  Col = 5; New = "554×9=4986" },

    @{ Row = 10; Col = 1; New = "775×4=3100" },
    @{ Row = 10; Col = 2; New = "621×4=2484" },
    @{ Row = 10; Col = 3; New = "522×4=2088" },
    @{ Row = 10; Col = 4; New = "486×3=1458" },
    @{ Row = 10; Col = 5; New = "669×9=6021" },

    @{ Row = 15; Col = 1; New = "935×7=6545" },
    @{ Row = 15; Col = 2; New = "531×4=2124" },
    @{ Row = 15; Col = 3; New = "292×7=2044" },
    @{ Row = 15; Col = 4; New = "527×4=2108" },
    @{ Row = 15; Col = 5; New = "589×5=2945" },

    @{ Row = 20; Col = 1; New = "658×6=3948" },
    @{ Row = 20; Col = 2; New = "124×6=744" },
    @{ Row = 20; Col = 3; New = "104×7=728" },
    @{ Row = 20; Col = 4; New = "678×4=2712" },
    @{ Row = 20; Col = 5; New = "847×5=4235" }
)

foreach ($item in $replacements) {
    $cell = $t.Cell($item.Row, $item.Col)
    $cellRange = $cell.Range
    # Exclude the trailing end-of-cell mark (the last char of Range.Text)
    # so only the visible content is replaced, preserving formatting.
    $contentRange = $d.Range($cellRange.Start, $cellRange.End - 1)
    $contentRange.Text = $item.New
}
